$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly record data between row 2 and row 5
# (Fecha, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)

$row2Date = $ws.Range("D2").Value2
$row5Date = $ws.Range("D5").Value2
$ws.Range("D2").Value2 = $row5Date
$ws.Range("D5").Value2 = $row2Date

$row2J = $ws.Range("J2").Value2
$row5J = $ws.Range("J5").Value2
$ws.Range("J2").Value2 = $row5J
$ws.Range("J5").Value2 = $row2J

$row2K = $ws.Range("K2").Value2
$row5K = $ws.Range("K5").Value2
$ws.Range("K2").Value2 = $row5K
$ws.Range("K5").Value2 = $row2K

$row2L = $ws.Range("L2").Value2
$row5L = $ws.Range("L5").Value2
$ws.Range("L2").Value2 = $row5L
$ws.Range("L5").Value2 = $row2L

$row2M = $ws.Range("M2").Value2
$row5M = $ws.Range("M5").Value2
$ws.Range("M2").Value2 = $row5M
$ws.Range("M5").Value2 = $row2M

$row2P = $ws.Range("P2").Value2
$row5P = $ws.Range("P5").Value2
$ws.Range("P2").Value2 = $row5P
$ws.Range("P5").Value2 = $row2P
